$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.47'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '0.91%'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '-2.25%'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.152'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '0.18%'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05777'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '1.84%'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.614'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '1.42%'
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.170'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '4.94%'
$ws.Range("B8").Value = 'MXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8564'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '2.00%'
$ws.Range("B9").Value = 'FTXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.8534'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '-1.28%'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1365'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2.20%'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07068'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '2.28%'
$ws.Range("B12").Value = 'BitrueCoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03050'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '6.51%'
$ws.Range("B13").Value = 'BitMartToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.09369'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '-0.13%'
$ws.Range("B14").Value = 'BitForexToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.001531'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '0.60%'
$ws.Range("B15").Value = 'One'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0006019'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '0.53%'
$ws.Range("B16").Value = 'TigerCash'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.006026'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '-0.49%'
$ws.Range("B17").Value = 'LEO'
$ws.Range("C17").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.485'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '-0.63%'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.161'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '-3.37%'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.03308'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '1.67%'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1284'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '-0.89%'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.318'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '-8.69%'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.04128'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '-0.85%'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '1.93%'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '1.25%'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004132'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '-4.31%'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '2.59%'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '3.43%'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03726'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '0.33%'
$ws.Range("B41").Value = 'KickToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.005890'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '10.56%'
$ws.Range("B42").Value = 'BKEXToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1069'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '1.09%'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002200'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '-4.77%'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.009151'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '-5.48%'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005294'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '3.61%'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '0.04%'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.05799'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '-41.97%'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002173'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '-20.01%'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.00002100'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '0.04%'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0002000'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '0.04%'
